$d = $word.ActiveDocument
$quoteL = [char]0x201C
$quoteR = [char]0x201D

$rPr = '<w:rPr><w:b/><w:bCs/><w:sz w:val="20"/></w:rPr>'
$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) "First-Class Mail" paragraph: split the single run into four runs and
#    wrap the mailing address in a Jinja {% if %}...{% endif %} block.
# ---------------------------------------------------------------------------
$oldMail = "{{output_checkbox(notice_type == " + $quoteL + "mail" + $quoteR + ")}} First-Class Mail, to this mailing address: {{ service_address.on_one_line() }}"

$rngMail = $d.Content
$foundMail = $rngMail.Find.Execute($oldMail, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundMail) {
    throw "Could not find the First-Class Mail paragraph text"
}

$pPrMail = '<w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:before="13" w:line="276" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/></w:rPr></w:pPr>'

$mailRun1 = '<w:r w:rsidRPr="005E072D">' + $rPr + '<w:t xml:space="preserve">{{output_checkbox(notice_type == ' + $quoteL + 'mail' + $quoteR + ')}} First-Class Mail, to this mailing address: </w:t></w:r>'
$mailRun2 = '<w:r>' + $rPr + '<w:t>{% if notice_type ==' + $quoteR + 'mail' + $quoteR + ' %}</w:t></w:r>'
$mailRun3 = '<w:r>' + $rPr + '<w:t>{{ service_address.on_one_line() }}</w:t></w:r>'
$mailRun4 = '<w:r>' + $rPr + '<w:t>{% endif %}</w:t></w:r>'

$newMailP = '<w:p w14:paraId="2B69B7B8" w14:textId="77777777" w:rsidR="005E072D" w:rsidRPr="005E072D" w:rsidRDefault="005E072D" w:rsidP="005E072D">' + $pPrMail + $mailRun1 + $mailRun2 + $mailRun3 + $mailRun4 + '</w:p>'

$rngMail.InsertXML($pkgHeader + $newMailP + $pkgFooter)

# ---------------------------------------------------------------------------
# 2) "E-Mail" paragraph: same transformation for the e-mail address.
# ---------------------------------------------------------------------------
$oldEmail = "{{output_checkbox(notice_type == " + $quoteL + "email" + $quoteR + ")}} E-Mail, to this e-mail address: {{ service_email_address }}"

$rngEmail = $d.Content
$foundEmail = $rngEmail.Find.Execute($oldEmail, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundEmail) {
    throw "Could not find the E-Mail paragraph text"
}

$pPrEmail = '<w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:before="13" w:line="276" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/></w:rPr></w:pPr>'

$emailRun1 = '<w:r w:rsidRPr="005E072D">' + $rPr + '<w:t xml:space="preserve">{{output_checkbox(notice_type == ' + $quoteL + 'email' + $quoteR + ')}} E-Mail, to this e-mail address: </w:t></w:r>'
$emailRun2 = '<w:r>' + $rPr + '<w:t>{% if notice_type == ' + $quoteL + 'email' + $quoteR + ' %}</w:t></w:r>'
$emailRun3 = '<w:r>' + $rPr + '<w:t>{{ service_email_address }}</w:t></w:r>'
$emailRun4 = '<w:r>' + $rPr + '<w:t>{% endif %}</w:t></w:r>'

$newEmailP = '<w:p w14:paraId="1544B4C0" w14:textId="2453123D" w:rsidR="00FE41B5" w:rsidRPr="005348D7" w:rsidRDefault="005E072D" w:rsidP="005E072D">' + $pPrEmail + $emailRun1 + $emailRun2 + $emailRun3 + $emailRun4 + '</w:p>'

$rngEmail.InsertXML($pkgHeader + $newEmailP + $pkgFooter)

Write-Output "OK"
